$d = $word.ActiveDocument

# The "COMPETENCES TECHNIQUES" section contains 8 single-line paragraphs that
# get reordered by this change (same 8 lines, new sequence). Locate the
# section by its heading text, then rewrite the 8 following paragraphs with
# the new order while leaving their paragraph formatting untouched.

$newOrder = @(
    "Langages : scala, python, matlab, c, c++",
    "Visualisation : etl, tableau",
    "MLOps : aws, databases, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit",
    "Web : api",
    "Autres : elasticsearch, dubai, posgresql, handy and pragmatic, consolidating, oil and gas companies, test, london",
    "Maths : algorithms",
    "ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn",
    "Bases de données : SQL, MongoDB, Neo4j, Redis"
)

$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "COMPETENCES TECHNIQUES") {
        $headingIndex = $i
        break
    }
}

for ($j = 0; $j -lt $newOrder.Length; $j++) {
    $p = $d.Paragraphs.Item($headingIndex + 1 + $j)
    $p.Range.Text = $newOrder[$j]
}
